$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells keep their original text representation
# (avoids Excel auto-converting numeric-looking strings to numbers,
# which would strip leading/trailing zeros and dot-grouping).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.775.03"
$ws.Range("E2").Value = "  +4.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.770.81"
$ws.Range("E3").Value = "  +4.53%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "343.84"
$ws.Range("E5").Value = "  +4.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "114.94"
$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.548"
$ws.Range("E7").Value = "  +4.42%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("E9").Value = "  +4.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.42"
$ws.Range("E10").Value = "  +6.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0853"
$ws.Range("E11").Value = "  +4.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.95"
$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.131"
$ws.Range("E13").Value = "  +1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.63"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.205.53"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.774.08"
$ws.Range("E16").Value = "  +4.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.878"
$ws.Range("E17").Value = "  +1.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.656.39"
$ws.Range("E18").Value = "  +3.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.21"
$ws.Range("E19").Value = "  +9.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.02"
$ws.Range("E20").Value = "  +4.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.19"
$ws.Range("E21").Value = "  -2.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0976"
$ws.Range("E22").Value = "  +2.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "271.80"
$ws.Range("E23").Value = "  +0.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.86"
$ws.Range("E24").Value = "  +0.93%  "

$ws.Range("E25").Value = "  +7.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.45"
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.19"
$ws.Range("E28").Value = "  -0.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  +1.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.139"
$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.35"
$ws.Range("E31").Value = "  -2.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.04"
$ws.Range("E32").Value = "  +0.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.69"
$ws.Range("E33").Value = "  +3.34%  "

$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.90"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.09"
$ws.Range("E37").Value = "  +2.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.93"
$ws.Range("E38").Value = "  -0.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.21"
$ws.Range("E39").Value = "  +2.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0381"
$ws.Range("E40").Value = "  +10.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.67"
$ws.Range("E41").Value = "  +24.95%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  +3.38%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.25"
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.33"
$ws.Range("E44").Value = "  +1.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.21"
$ws.Range("E45").Value = "  -4.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.065.07"
$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("E48").Value = "  +1.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.52"
$ws.Range("E49").Value = "  +4.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.890"
$ws.Range("E50").Value = "  +12.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.86"
$ws.Range("E51").Value = "  -1.42%  "
